$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule entries appended below the existing table (rows 53-56),
# mirroring the same per-row layout used throughout the sheet (row 52 is a
# representative fully-styled row: column A = name/style-11, columns B-J the
# plain "cell" style-12, with the style-13 "text" variant on whichever column
# actually carries a label for that week).

# --- Row 53 : Nargiz Ahmadova (007) / GROUP 18 in column J ---
$ws.Range("A52:J52").Copy()
$ws.Range("A53:J53").PasteSpecial(-4122)
$ws.Rows.Item(53).RowHeight = 13.65
$ws.Cells.Item(53,1).Value = "Nargiz Ahmadova (007)"
$ws.Cells.Item(53,10).Value = "GROUP 18"
$ws.Cells.Item(53,10).NumberFormat = "@"

# --- Row 54 : Altaf Hussain / GROUP 16 in column I ---
$ws.Range("A52:J52").Copy()
$ws.Range("A54:J54").PasteSpecial(-4122)
$ws.Rows.Item(54).RowHeight = 13.65
$ws.Cells.Item(54,1).Value = "Altaf Hussain"
$ws.Cells.Item(54,9).Value = "GROUP 16"
$ws.Cells.Item(54,9).NumberFormat = "@"

# --- Row 55 : Valentyn Khmarskyi / GROUP 5 in column D ---
$ws.Range("A52:J52").Copy()
$ws.Range("A55:J55").PasteSpecial(-4122)
$ws.Rows.Item(55).RowHeight = 13.65
$ws.Cells.Item(55,1).Value = "Valentyn Khmarskyi"
$ws.Cells.Item(55,4).Value = "GROUP 5"
$ws.Cells.Item(55,4).NumberFormat = "@"

# --- Row 56 : Ozodbek Ozodov / GROUP 13 in column H ---
$ws.Range("A52:J52").Copy()
$ws.Range("A56:J56").PasteSpecial(-4122)
$ws.Rows.Item(56).RowHeight = 13.65
$ws.Cells.Item(56,1).Value = "Ozodbek Ozodov"
$ws.Cells.Item(56,8).Value = "GROUP 13"
$ws.Cells.Item(56,8).NumberFormat = "@"

$excel.CutCopyMode = $false
